$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (player, position, team) replacing the old row order in A2:C19.
$rows = @(
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("D'Angelo Russell", "PG", "Los Angeles Lakers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Jrue Holiday", "PG,SG", "Boston Celtics"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
